$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.114.28"
$ws.Range("E2").Value = "  +4.29%  "

$ws.Range("D3").Value = "1.906.02"
$ws.Range("E3").Value = "  +5.26%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9994"
$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "252.99"
$ws.Range("E5").Value = "  +2.01%  "

$ws.Range("E6").Value = "  -0.05%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5083"
$ws.Range("E7").Value = "  +2.70%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "45.21"
$ws.Range("E8").Value = "  +4.50%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3037"
$ws.Range("E9").Value = "  +9.09%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06818"
$ws.Range("E10").Value = "  +6.20%  "

$ws.Range("D11").Value = "1.907.93"
$ws.Range("E11").Value = "  +5.44%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "17.30"
$ws.Range("E12").Value = "  +2.94%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07321"
$ws.Range("E13").Value = "  +3.24%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6914"
$ws.Range("E14").Value = "  +6.99%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "86.93"
$ws.Range("E15").Value = "  +3.22%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.915"
$ws.Range("E16").Value = "  +4.85%  "

$ws.Range("D17").Value = "30.117.17"
$ws.Range("E17").Value = "  +4.18%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008284"
$ws.Range("E18").Value = "  +12.83%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9984"
$ws.Range("E19").Value = "  -0.07%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.06"
$ws.Range("E20").Value = "  +6.47%  "

$ws.Range("D21").Value = "2.151.95"
$ws.Range("E21").Value = "  +5.14%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9992"
$ws.Range("E22").Value = "  -0.02%  "

$ws.Range("E23").Value = "  +5.25%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.739"
$ws.Range("E24").Value = "  +7.07%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.297"
$ws.Range("E25").Value = "  +5.01%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "147.61"
$ws.Range("E26").Value = "  +3.71%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "134.72"
$ws.Range("E27").Value = "  +4.50%  "

$ws.Range("E28").Value = "  +4.39%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.004"
$ws.Range("E29").Value = "  +5.53%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.399"
$ws.Range("E30").Value = "  -1.06%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.282"
$ws.Range("E31").Value = "  +3.17%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08860"
$ws.Range("E32").Value = "  +6.04%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.005"
$ws.Range("E33").Value = "  +4.92%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05062"
$ws.Range("E34").Value = "  +1.87%  "

$ws.Range("E35").Value = "  +4.19%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7231"
$ws.Range("E36").Value = "  +7.24%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.689"
$ws.Range("E37").Value = "  +0.16%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.822"
$ws.Range("E38").Value = "  +2.65%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.274"
$ws.Range("E39").Value = "  -2.52%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9630"
$ws.Range("E40").Value = "  +0.99%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.01691"
$ws.Range("E41").Value = "  +6.07%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.167"
$ws.Range("E42").Value = "  +0.24%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4313"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "104.94"
$ws.Range("E44").Value = "  +5.15%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9989"
$ws.Range("E45").Value = "  -0.01%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.661"
$ws.Range("E46").Value = "  +7.06%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1280"
$ws.Range("E47").Value = "  +4.89%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05750"
$ws.Range("E48").Value = "  +4.24%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "33.25"
$ws.Range("E49").Value = "  +4.80%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.409"
$ws.Range("E50").Value = "  +2.78%  "

$ws.Range("E51").Value = "  +5.25%  "
